$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (jeremy o)
$ws.Range("C2").Value = 87
$ws.Range("D2").Value = 75
$ws.Range("E2").Value = 90
$ws.Range("F2").Value = 89

# Row 3 (jeremy c)
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 90
$ws.Range("F3").Value = 93

# Row 4 (remy)
$ws.Range("C4").Value = 81
$ws.Range("D4").Value = 88
$ws.Range("E4").Value = 87
$ws.Range("F4").Value = 89

# Row 5 (cristina)
$ws.Range("C5").Value = 80
$ws.Range("D5").Value = 86
$ws.Range("E5").Value = 88
$ws.Range("F5").Value = 91
